$d = $word.ActiveDocument

# Explicitly set "page break before" to False on every paragraph in the
# document body (this materializes <w:pageBreakBefore w:val="0"/> inside
# each paragraph's <w:pPr>, matching the authoring tool's output).
$d.Paragraphs.PageBreakBefore = $false

# Also set it explicitly (to False) on the built-in heading / title
# styles used in the document, which previously had no pageBreakBefore
# setting at all.
$styleNames = @("Heading1", "Heading2", "Heading3", "Heading4", "Heading5", "Heading6", "Title", "Subtitle")
foreach ($name in $styleNames) {
    $style = $d.Styles($name)
    $style.ParagraphFormat.PageBreakBefore = $false
}

Write-Output "pageBreakBefore normalized on paragraphs and heading/title styles"
